$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Formula = "=" + "CONCATENATE(""INSERT INTO ``clinic`` (``"",A1,""``,``"",B1,""``,``"",C1,""``,``"",D1,""``,``"",E1,""``,``"",F1,""``,``"",G1,""``,``"",H1,""``,``"",I1,""``,``"",J1,""``,``"",K1,""``,``"",L1,""``,``"",M1,""``,``"",N1,""``,``"",O1,""``,``"",P1,""``,``"",Q1,""``) VALUES ("")"
$ws.Range("AM1").Value = "SQL Command"
$ws.Range("AM1").Font.Bold = $true

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 0

$ws.Range("S2").Formula = "=" + "IF(A2<>"""",CONCATENATE("""""""",A2,"""""",""),""NULL,"")"
$ws.Range("T2").Formula = "=" + "IF(B2<>"""",CONCATENATE("""""""",B2,"""""",""),""NULL,"")"
$ws.Range("U2").Formula = "=" + "IF(C2<>"""",CONCATENATE("""""""",C2,"""""",""),""NULL,"")"
$ws.Range("V2").Formula = "=" + "IF(D2<>"""",CONCATENATE("""""""",D2,"""""",""),""NULL,"")"
$ws.Range("W2").Formula = "=" + "IF(E2<>"""",CONCATENATE("""""""",E2,"""""",""),""NULL,"")"
$ws.Range("X2").Formula = "=" + "IF(F2<>"""",CONCATENATE("""""""",F2,"""""",""),""NULL,"")"
$ws.Range("Y2").Formula = "=" + "IF(G2<>"""",CONCATENATE("""""""",G2,"""""",""),""NULL,"")"
$ws.Range("Z2").Formula = "=" + "IF(H2<>"""",CONCATENATE("""""""",H2,"""""",""),""NULL,"")"
$ws.Range("AA2").Formula = "=" + "IF(I2<>"""",CONCATENATE("""""""",I2,"""""",""),""NULL,"")"
$ws.Range("AB2").Formula = "=" + "IF(J2<>"""",CONCATENATE("""""""",J2,"""""",""),""NULL,"")"
$ws.Range("AC2").Formula = "=" + "IF(K2<>"""",CONCATENATE("""""""",K2,"""""",""),""NULL,"")"
$ws.Range("AD2").Formula = "=" + "IF(L2<>"""",CONCATENATE("""""""",L2,"""""",""),""NULL,"")"
$ws.Range("AE2").Formula = "=" + "IF(M2<>"""",CONCATENATE("""""""",M2,"""""",""),""NULL,"")"
$ws.Range("AF2").Formula = "=" + "IF(N2<>"""",CONCATENATE("""""""",N2,"""""",""),""NULL,"")"
$ws.Range("AG2").Formula = "=" + "IF(O2<>"""",CONCATENATE("""""""",O2,"""""",""),""NULL,"")"
$ws.Range("AH2").Formula = "=" + "IF(P2<>"""",CONCATENATE(P2,"",""),""NULL,"")"
$ws.Range("AI2").Formula = "=" + "IF(Q2<>"""",CONCATENATE(Q2),""NULL"")"
$ws.Range("AK2").Formula = "=" + "CONCATENATE(S2,T2,U2,V2,W2,X2,Y2,Z2,AA2,AB2,AC2,AD2,AE2,AF2,AG2,AH2,AI2,"");"")"
$ws.Range("AM2").Formula = "=" + "CONCATENATE(`$S`$1,AK2)"

$ws.Range("S3").Formula = "=" + "IF(A3<>"""",CONCATENATE("""""""",A3,"""""",""),""NULL,"")"
$ws.Range("T3").Formula = "=" + "IF(B3<>"""",CONCATENATE("""""""",B3,"""""",""),""NULL,"")"
$ws.Range("U3").Formula = "=" + "IF(C3<>"""",CONCATENATE("""""""",C3,"""""",""),""NULL,"")"
$ws.Range("V3").Formula = "=" + "IF(D3<>"""",CONCATENATE("""""""",D3,"""""",""),""NULL,"")"
$ws.Range("W3").Formula = "=" + "IF(E3<>"""",CONCATENATE("""""""",E3,"""""",""),""NULL,"")"
$ws.Range("X3").Formula = "=" + "IF(F3<>"""",CONCATENATE("""""""",F3,"""""",""),""NULL,"")"
$ws.Range("Y3").Formula = "=" + "IF(G3<>"""",CONCATENATE("""""""",G3,"""""",""),""NULL,"")"
$ws.Range("Z3").Formula = "=" + "IF(H3<>"""",CONCATENATE("""""""",H3,"""""",""),""NULL,"")"
$ws.Range("AA3").Formula = "=" + "IF(I3<>"""",CONCATENATE("""""""",I3,"""""",""),""NULL,"")"
$ws.Range("AB3").Formula = "=" + "IF(J3<>"""",CONCATENATE("""""""",J3,"""""",""),""NULL,"")"
$ws.Range("AC3").Formula = "=" + "IF(K3<>"""",CONCATENATE("""""""",K3,"""""",""),""NULL,"")"
$ws.Range("AD3").Formula = "=" + "IF(L3<>"""",CONCATENATE("""""""",L3,"""""",""),""NULL,"")"
$ws.Range("AE3").Formula = "=" + "IF(M3<>"""",CONCATENATE("""""""",M3,"""""",""),""NULL,"")"
$ws.Range("AF3").Formula = "=" + "IF(N3<>"""",CONCATENATE("""""""",N3,"""""",""),""NULL,"")"
$ws.Range("AG3").Formula = "=" + "IF(O3<>"""",CONCATENATE("""""""",O3,"""""",""),""NULL,"")"
$ws.Range("AH3").Formula = "=" + "IF(P3<>"""",CONCATENATE(P3,"",""),""NULL,"")"
$ws.Range("AI3").Formula = "=" + "IF(Q3<>"""",CONCATENATE(Q3),""NULL"")"
$ws.Range("AK3").Formula = "=" + "CONCATENATE(S3,T3,U3,V3,W3,X3,Y3,Z3,AA3,AB3,AC3,AD3,AE3,AF3,AG3,AH3,AI3,"");"")"
$ws.Range("AM3").Formula = "=" + "CONCATENATE(`$S`$1,AK3)"

$ws.Range("S4").Formula = "=" + "IF(A4<>"""",CONCATENATE("""""""",A4,"""""",""),""NULL,"")"
$ws.Range("T4").Formula = "=" + "IF(B4<>"""",CONCATENATE("""""""",B4,"""""",""),""NULL,"")"
$ws.Range("U4").Formula = "=" + "IF(C4<>"""",CONCATENATE("""""""",C4,"""""",""),""NULL,"")"
$ws.Range("V4").Formula = "=" + "IF(D4<>"""",CONCATENATE("""""""",D4,"""""",""),""NULL,"")"
$ws.Range("W4").Formula = "=" + "IF(E4<>"""",CONCATENATE("""""""",E4,"""""",""),""NULL,"")"
$ws.Range("X4").Formula = "=" + "IF(F4<>"""",CONCATENATE("""""""",F4,"""""",""),""NULL,"")"
$ws.Range("Y4").Formula = "=" + "IF(G4<>"""",CONCATENATE("""""""",G4,"""""",""),""NULL,"")"
$ws.Range("Z4").Formula = "=" + "IF(H4<>"""",CONCATENATE("""""""",H4,"""""",""),""NULL,"")"
$ws.Range("AA4").Formula = "=" + "IF(I4<>"""",CONCATENATE("""""""",I4,"""""",""),""NULL,"")"
$ws.Range("AB4").Formula = "=" + "IF(J4<>"""",CONCATENATE("""""""",J4,"""""",""),""NULL,"")"
$ws.Range("AC4").Formula = "=" + "IF(K4<>"""",CONCATENATE("""""""",K4,"""""",""),""NULL,"")"
$ws.Range("AD4").Formula = "=" + "IF(L4<>"""",CONCATENATE("""""""",L4,"""""",""),""NULL,"")"
$ws.Range("AE4").Formula = "=" + "IF(M4<>"""",CONCATENATE("""""""",M4,"""""",""),""NULL,"")"
$ws.Range("AF4").Formula = "=" + "IF(N4<>"""",CONCATENATE("""""""",N4,"""""",""),""NULL,"")"
$ws.Range("AG4").Formula = "=" + "IF(O4<>"""",CONCATENATE("""""""",O4,"""""",""),""NULL,"")"
$ws.Range("AH4").Formula = "=" + "IF(P4<>"""",CONCATENATE(P4,"",""),""NULL,"")"
$ws.Range("AI4").Formula = "=" + "IF(Q4<>"""",CONCATENATE(Q4),""NULL"")"
$ws.Range("AK4").Formula = "=" + "CONCATENATE(S4,T4,U4,V4,W4,X4,Y4,Z4,AA4,AB4,AC4,AD4,AE4,AF4,AG4,AH4,AI4,"");"")"
$ws.Range("AM4").Formula = "=" + "CONCATENATE(`$S`$1,AK4)"

$ws.Range("S1:AL1").EntireColumn.AutoFit()

$ws.Range("AF1").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 32
$win.ScrollRow = 1
$ws.Range("AK2").Select()
